$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 135, shifting existing rows 135:177 down to 136:178
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new data record
$ws.Cells.Item(135, 1).Value2 = 10
$ws.Cells.Item(135, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(135, 3).Value2 = "La Araucanía"
$ws.Cells.Item(135, 4).Value2 = 44985
$ws.Cells.Item(135, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(135, 5).Value2 = 9
$ws.Cells.Item(135, 6).Value2 = 100114002
$ws.Cells.Item(135, 7).Value2 = "Camote"
$ws.Cells.Item(135, 8).Value2 = "Sin especificar"
$ws.Cells.Item(135, 9).Value2 = "Primera"
$ws.Cells.Item(135, 10).Value2 = 30
$ws.Cells.Item(135, 11).Value2 = 26000
$ws.Cells.Item(135, 12).Value2 = 26000
$ws.Cells.Item(135, 13).Value2 = 26000
$ws.Cells.Item(135, 14).Value2 = "$/malla 20 kilos"
$ws.Cells.Item(135, 15).Value2 = "Perú"
$ws.Cells.Item(135, 16).Value2 = 1300
$ws.Cells.Item(135, 17).Value2 = 20
$ws.Cells.Item(135, 18).Value2 = "Hortaliza"
